# "Generate Report for Handback"
# The handback tooling re-ran and produced fresh Handoff/Handback timestamps
# for the file f1b02328-bd5f-43ba-acd3-828da2075d44.md (row 6 on each sheet).
# Everything else in the report is unchanged.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G6").Value = "2016-10-27 02:14:58"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H6").Value = "2016-10-27 02:14:43"
$zhcn.Range("K6").Value = "2016-10-27 02:15:33"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H6").Value = "2016-10-27 02:14:58"
$dede.Range("K6").Value = "2016-10-27 02:15:51"
